# "Week Before Final PC Dump"
# Updates the cardinfo sheet: correct a mis-read byte value, document a new
# application flag byte, fill in a few of the sequential index cells that
# had been left blank, and move the selection to where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cardinfo")

# Document what the application-presence byte means.
$ws.Range("T3").Value = "Gives App precense and read write information for application"

# The byte that had been transcribed as 0x80 in row 3 ("Application And
# Validitiy Map") was actually 0x88 - fix the reading.
$ws.Range("N3").Value = " 0x88"

# Fill in the sequential application/data indices that were still blank.
$ws.Range("R8").Value = 3
$ws.Range("R10").Value = 4
$ws.Range("R12").Value = 5
$ws.Range("R14").Value = 6

# Leave the selection where the review work stopped.
$ws.Range("C8").Select()
